$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P1").Value = 0.74478269137660513
$ws.Range("Y1").Value = 0.86871639297336856
$ws.Range("N2").Value = 0.83972540022420783
$ws.Range("O2").Value = 0.80564809957498151
$ws.Range("AC2").Value = 0.78951482443259291
$ws.Range("G3").Value = 0.92478789733156308
$ws.Range("S3").Value = 0.93920361829316068
$ws.Range("AG3").Value = 0.93085247825338302
$ws.Range("AH3").Value = 0.92861389257124649
$ws.Range("AS3").Value = 0.91836400906875904
$ws.Range("S4").Value = 0.66788274186792551
$ws.Range("R5").Value = 0.85447599597542889
$ws.Range("AR5").Value = 0.76545512560411211
$ws.Range("BM5").Value = 0.88643104882169454
$ws.Range("AH7").Value = 0.77459696364303432
$ws.Range("BK7").Value = 0.85729560817937234
$ws.Range("X8").Value = 0.59016676190145689
$ws.Range("BI8").Value = 0.86813672639724104
$ws.Range("Y9").Value = 0.94441532911423742
$ws.Range("AJ9").Value = 0.85089074256132369
$ws.Range("AU9").Value = 0.79545033863998316
$ws.Range("BN9").Value = 0.74613065852789573
$ws.Range("D10").Value = 0.89268189275821297
$ws.Range("H10").Value = 0.96693991289172332
$ws.Range("I11").Value = 0.83618734579241982
$ws.Range("AS11").Value = 0.97239092022269968
$ws.Range("D12").Value = 0.85463236553537891
$ws.Range("P12").Value = 0.89259545979070309
$ws.Range("AS12").Value = 0.82706457664623156
$ws.Range("BL12").Value = 0.77001372166093285
$ws.Range("N13").Value = 0.67292709483262736
$ws.Range("AC13").Value = 0.94874498848181821
$ws.Range("AQ13").Value = 0.99395590000027312
$ws.Range("AT13").Value = 0.81804034740638376
$ws.Range("Y14").Value = 0.62231787816457684
$ws.Range("BI14").Value = 0.93079282164802091
$ws.Range("W15").Value = 0.73444831834121915
$ws.Range("AJ15").Value = 0.81039377663609791
$ws.Range("AT15").Value = 0.84256304719119224
$ws.Range("J16").Value = 0.97924161697505552
$ws.Range("R16").Value = 0.7337316873449331
$ws.Range("BO16").Value = 0.59993545979836627
$ws.Range("X17").Value = 0.95089129900357849
$ws.Range("J18").Value = 0.86915420854519088
$ws.Range("AE19").Value = 0.79034127011181909
$ws.Range("V20").Value = 0.89126252765836378
$ws.Range("D21").Value = 0.96878270892933782
$ws.Range("S21").Value = 0.879960922166515
$ws.Range("BO21").Value = 0.73882635615628911
$ws.Range("M22").Value = 0.88392469614728042
$ws.Range("X23").Value = 0.67827910063526509
$ws.Range("Z24").Value = 0.86313603708039821
$ws.Range("AV24").Value = 0.98363131115038716
$ws.Range("AZ25").Value = 0.9787762324105671
$ws.Range("AU26").Value = 0.83147878010870924
$ws.Range("BJ26").Value = 0.97610825430369474
$ws.Range("BP26").Value = 0.63281492559959407
$ws.Range("A27").Value = 0.94854190486094303
$ws.Range("BL27").Value = 0.76796261143366307
$ws.Range("S28").Value = 0.89615717055491495
$ws.Range("AB29").Value = 0.89149092779722683
$ws.Range("N30").Value = 0.67514617081048367
$ws.Range("AC30").Value = 0.99765280108173826
$ws.Range("AH30").Value = 0.86464154037652008
$ws.Range("AM31").Value = 0.83991436895465654
$ws.Range("BF31").Value = 0.63993461310252553
$ws.Range("BM32").Value = 0.85342062184121881
$ws.Range("T33").Value = 0.71708980888129514
$ws.Range("Y33").Value = 0.71526955429305106
$ws.Range("AF33").Value = 0.75393654293473533
$ws.Range("AR34").Value = 0.93171542695968723
$ws.Range("BH34").Value = 0.58792592146814227
$ws.Range("L35").Value = 0.99739577760881193
$ws.Range("T35").Value = 0.69779672445025087
$ws.Range("BJ36").Value = 0.6843950199938289
$ws.Range("F37").Value = 0.9711117831239845
$ws.Range("AN37").Value = 0.87574537418556075
$ws.Range("BA37").Value = 0.99801946338934533
$ws.Range("AG38").Value = 0.92834807067772496
$ws.Range("A39").Value = 0.60107452918138438
$ws.Range("AK39").Value = 0.79120659685436057
$ws.Range("AL39").Value = 0.66522598976421254
$ws.Range("AN39").Value = 0.87366074151703477
$ws.Range("AO39").Value = 0.81509672068914896
$ws.Range("I40").Value = 0.77267634474107116
$ws.Range("AX40").Value = 0.66052038420496717
$ws.Range("Q41").Value = 0.94802954527822259
$ws.Range("AQ41").Value = 0.88632162217565302
$ws.Range("BC41").Value = 0.81453071117479525
$ws.Range("BG41").Value = 0.78479158152399897
$ws.Range("U42").Value = 0.97796992183375098
$ws.Range("W43").Value = 0.58296260786593945
$ws.Range("AS43").Value = 0.91702492581376105
$ws.Range("AA45").Value = 0.94781949288991563
$ws.Range("AR46").Value = 0.74064878067758122
$ws.Range("F48").Value = 0.76026722134106661
$ws.Range("BA48").Value = 0.95839740321529743
$ws.Range("BE48").Value = 0.68336967303751339
$ws.Range("N49").Value = 0.91154621310926631
$ws.Range("V49").Value = 0.96408305405429551
$ws.Range("Q50").Value = 0.98290722398247476
$ws.Range("AZ50").Value = 0.83370413866631898
$ws.Range("BC50").Value = 0.91070683473883496
$ws.Range("V51").Value = 0.88265775758508125
$ws.Range("BB51").Value = 0.86078421507079039
$ws.Range("BA52").Value = 0.7468601294049626
$ws.Range("P53").Value = 0.71396101663171385
$ws.Range("BH53").Value = 0.94068914167985451
$ws.Range("M55").Value = 0.96877647408799472
$ws.Range("U55").Value = 0.94203081603756655
$ws.Range("BB55").Value = 0.88510903671231866
$ws.Range("BD55").Value = 0.76079642726255525
$ws.Range("R56").Value = 0.94787750345834887
$ws.Range("N57").Value = 0.76059711441619227
$ws.Range("AP57").Value = 0.90839155073015121
$ws.Range("E58").Value = 0.7716565369164381
$ws.Range("AF59").Value = 0.86524964250888936
$ws.Range("AZ59").Value = 0.97796523741885055
$ws.Range("AS60").Value = 0.95203196312068084
$ws.Range("AY60").Value = 0.84831732153390838
$ws.Range("R62").Value = 0.8092378187876681
$ws.Range("BL62").Value = 0.72756893219003027
$ws.Range("AW63").Value = 0.96317594096113712
$ws.Range("BJ63").Value = 0.94775297321220964
$ws.Range("BL63").Value = 0.8804591151920188
$ws.Range("N64").Value = 0.81941184073499662
$ws.Range("R64").Value = 0.90897728747908091
$ws.Range("AU64").Value = 0.67200273967823221
$ws.Range("K66").Value = 0.87878364683661303
$ws.Range("Z66").Value = 0.99034947583797361
$ws.Range("AF66").Value = 0.91244602014096943
$ws.Range("AG66").Value = 0.86524914596669744
$ws.Range("BC66").Value = 0.94428024478112293
$ws.Range("BM66").Value = 0.96299015604760996
$ws.Range("A68").Value = 0.78950625827217447
$ws.Range("BI68").Value = 0.96944674764345229
